$d = $word.ActiveDocument

$replacements = @(
    @("283÷7=", "664÷2="),
    @("550÷4=", "168÷9="),
    @("377÷4=", "352÷3="),
    @("167÷2=", "633÷5="),
    @("312÷7=", "239÷8="),
    @("975÷3=", "517÷5="),
    @("894÷9=", "211÷6="),
    @("876÷9=", "668÷6="),
    @("581÷7=", "785÷3="),
    @("661÷8=", "726÷2="),
    @("307÷2=", "548÷8="),
    @("988÷3=", "573÷4="),
    @("575÷6=", "880÷6="),
    @("473÷4=", "170÷4="),
    @("173÷9=", "537÷3="),
    @("792÷8=", "433÷2="),
    @("140÷7=", "206÷7="),
    @("969÷4=", "850÷8="),
    @("784÷9=", "180÷6="),
    @("921÷8=", "409÷5="),
    @("230÷2=", "115÷3="),
    @("284÷9=", "404÷2="),
    @("308÷6=", "533÷7="),
    @("921÷2=", "125÷6="),
    @("178÷7=", "124÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
